# Github Auto Build at 2023-12-13 07:16
# Append the latest cost-log rows to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A258").Value = "2023-12-13 07:15:26"
$ws.Range("B258").Value = 0.0004

$ws.Range("A259").Value = "2023-12-13 07:15:45"
$ws.Range("B259").Value = 0.0006000000000000001

$ws.Range("A260").Value = "2023-12-13 07:16:03"
$ws.Range("B260").Value = 0.0008
